# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The worker/period detail table (rows 16-25) previously listed:
#   row16: CC 1002190783 SERGIO DAVID MORINELLY RAMOS  2507
#   row17: CC 1000130337 KAROL DAYAN VEGA GUTIERREZ    2507
#   row18: CC 1000130337 KAROL DAYAN VEGA GUTIERREZ    2506
#   row19: CC 1000130337 KAROL DAYAN VEGA GUTIERREZ    2505
#   row20: CC 1000130337 KAROL DAYAN VEGA GUTIERREZ    2504
#   row21: CC 1000130337 KAROL DAYAN VEGA GUTIERREZ    2503
#   row22: CC 1000130337 KAROL DAYAN VEGA GUTIERREZ    2502
#   row23: CC 1047483343 PATRICIA JUDITH RODRIGUEZ MONTES 2505
#   row24: CC 1047483343 PATRICIA JUDITH RODRIGUEZ MONTES 2504
#   row25: CC 1047483343 PATRICIA JUDITH RODRIGUEZ MONTES 2503
#
# The new database only keeps SERGIO DAVID MORINELLY RAMOS, with one extra
# period (2508) added as "parte 1" of the new account statement. So rows
# 17-24 (the KAROL / PATRICIA rows) are removed, and the remaining last row
# (which becomes the new row 17) is turned into the new SERGIO / 2508 entry.
# This also shifts the footer signature rows (old 30/31) up to 22/23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 8 rows (17-24) that belonged to the workers/periods no longer
# present in the updated database. This shifts everything below up, so the
# former row 25 becomes row 17, and the former rows 30/31 become 22/23.
$ws.Rows("17:24").Delete()

# The row that slid up into row 17 keeps the last-row border styling, but its
# data must now mirror row 16 (same worker) with the new period 2508.
$ws.Range("C17").Value = $ws.Range("C16").Value2
$ws.Range("D17").Value = $ws.Range("D16").Value2
$ws.Range("E17").Value = "2508"

# Refresh the summary totals for the smaller dataset.
$ws.Range("E11").Value = 113880
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 2
